# Auto-generated Excel COM-interop script that applies the scheduled-runner
# market-data refresh described in the commit diff.
# For each affected Leve row, the currentAveragePrice* (H/I/J) and the
# computed LevePrice*/LeveProfit* (K/L/M/N) columns are rewritten to the
# latest pulled values. A couple of rows additionally drop their N column
# (LeveProfitHQ) entirely because the HQ price is no longer available.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2026.3572
$ws.Range("I40").Value = 1843.5
$ws.Range("J40").Value = 2117.7856
$ws.Range("K40").Value = 1843.5
$ws.Range("L40").Value = 2117.7856
$ws.Range("M40").Value = -1668.5
$ws.Range("N40").Value = -2467.7856

$ws.Range("H93").Value = 28800
$ws.Range("J93").Value = 28800
$ws.Range("L93").Value = 28800
$ws.Range("N93").Value = -33792

$ws.Range("H129").Value = 1372618
$ws.Range("I129").Value = 263.69232
$ws.Range("J129").Value = 2646947
$ws.Range("K129").Value = 791.07696
$ws.Range("L129").Value = 7940841
$ws.Range("M129").Value = 4208.92304
$ws.Range("N129").Value = -7950841

$ws.Range("H132").Value = 3105.12
$ws.Range("I132").Value = 3498.9
$ws.Range("K132").Value = 10496.7
$ws.Range("M132").Value = -7966.700000000001

$ws.Range("H137").Value = 1109
$ws.Range("I137").Value = 1069.9
$ws.Range("K137").Value = 3209.7
$ws.Range("M137").Value = -659.7000000000003

$ws.Range("H138").Value = 1146.8315
$ws.Range("I138").Value = 614.7027
$ws.Range("J138").Value = 1486.2931
$ws.Range("K138").Value = 1844.1081
$ws.Range("L138").Value = 4458.879300000001
$ws.Range("M138").Value = 3295.8919
$ws.Range("N138").Value = -14738.8793

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5578556
$ws.Range("I32").Value = 6084497.5
$ws.Range("J32").Value = 13200
$ws.Range("K32").Value = 6084497.5
$ws.Range("L32").Value = 13200
$ws.Range("M32").Value = -6084210.5
$ws.Range("N32").Value = -13774

$ws.Range("H61").Value = 1165.2858
$ws.Range("I61").Value = 859.5
$ws.Range("K61").Value = 859.5
$ws.Range("M61").Value = -647.5

$ws.Range("H74").Value = 1412.4348
$ws.Range("I74").Value = 1412.4348
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1412.4348
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -538.4348
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1412.4348
$ws.Range("I77").Value = 1412.4348
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7062.174
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2694.174
$ws.Range("N77").ClearContents()

$ws.Range("H136").Value = 1165.2858
$ws.Range("I136").Value = 859.5
$ws.Range("K136").Value = 2578.5
$ws.Range("M136").Value = -28.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 25018660
$ws.Range("I20").Value = 31071.3
$ws.Range("J20").Value = 50006250
$ws.Range("K20").Value = 31071.3
$ws.Range("L20").Value = 50006250
$ws.Range("M20").Value = -30824.3
$ws.Range("N20").Value = -50006744

$ws.Range("H107").Value = 18864.25
$ws.Range("I107").Value = 2876.8333
$ws.Range("J107").Value = 66826.5
$ws.Range("K107").Value = 2876.8333
$ws.Range("L107").Value = 66826.5
$ws.Range("M107").Value = -956.8332999999998
$ws.Range("N107").Value = -70666.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 258
$ws.Range("I7").Value = 252.33333
$ws.Range("J7").Value = 292
$ws.Range("K7").Value = 252.33333
$ws.Range("L7").Value = 292
$ws.Range("M7").Value = -139.33333
$ws.Range("N7").Value = -518

$ws.Range("H31").Value = 2839.9106
$ws.Range("I31").Value = 2941.262
$ws.Range("J31").Value = 2535.8572
$ws.Range("K31").Value = 2941.262
$ws.Range("L31").Value = 2535.8572
$ws.Range("M31").Value = -2646.262
$ws.Range("N31").Value = -3125.8572

$ws.Range("H34").Value = 2839.9106
$ws.Range("I34").Value = 2941.262
$ws.Range("J34").Value = 2535.8572
$ws.Range("K34").Value = 2941.262
$ws.Range("L34").Value = 2535.8572
$ws.Range("M34").Value = -2739.262
$ws.Range("N34").Value = -2939.8572

$ws.Range("H62").Value = 4616.364
$ws.Range("I62").Value = 4850
$ws.Range("K62").Value = 4850
$ws.Range("M62").Value = -4226

$ws.Range("H65").Value = 4616.364
$ws.Range("I65").Value = 4850
$ws.Range("K65").Value = 24250
$ws.Range("M65").Value = -21130

$ws.Range("H99").Value = 2465.12
$ws.Range("I99").Value = 2315.7896
$ws.Range("J99").Value = 2938
$ws.Range("K99").Value = 2315.7896
$ws.Range("L99").Value = 2938
$ws.Range("M99").Value = -817.7896000000001
$ws.Range("N99").Value = -5934

$ws.Range("H122").Value = 1285.6666
$ws.Range("I122").Value = 1135.2593
$ws.Range("J122").Value = 1736.8889
$ws.Range("K122").Value = 3405.7779
$ws.Range("L122").Value = 5210.6667
$ws.Range("M122").Value = -955.7779
$ws.Range("N122").Value = -10110.6667

$ws.Range("H126").Value = 2465.12
$ws.Range("I126").Value = 2315.7896
$ws.Range("J126").Value = 2938
$ws.Range("K126").Value = 6947.3688
$ws.Range("L126").Value = 8814
$ws.Range("M126").Value = -4477.3688
$ws.Range("N126").Value = -13754

$ws.Range("H132").Value = 1449.6571
$ws.Range("I132").Value = 1038.75
$ws.Range("K132").Value = 3116.25
$ws.Range("M132").Value = -586.25

$ws.Range("H134").Value = 920.55884
$ws.Range("I134").Value = 774.2083
$ws.Range("J134").Value = 1271.8
$ws.Range("K134").Value = 2322.6249
$ws.Range("L134").Value = 3815.4
$ws.Range("M134").Value = 212.3751000000002
$ws.Range("N134").Value = -8885.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 42161.04
$ws.Range("I109").Value = 111825.11
$ws.Range("J109").Value = 2975
$ws.Range("K109").Value = 335475.33
$ws.Range("L109").Value = 8925
$ws.Range("M109").Value = -334435.33
$ws.Range("N109").Value = -11005

$ws.Range("H131").Value = 6408389
$ws.Range("I131").Value = 55101332
$ws.Range("J131").Value = 1422.921
$ws.Range("K131").Value = 165303996
$ws.Range("L131").Value = 4268.763
$ws.Range("M131").Value = -165298956
$ws.Range("N131").Value = -14348.763

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4789.778
$ws.Range("I70").Value = 4448.3335
$ws.Range("K70").Value = 4448.3335
$ws.Range("M70").Value = -4178.3335

$ws.Range("H73").Value = 4789.778
$ws.Range("I73").Value = 4448.3335
$ws.Range("K73").Value = 4448.3335
$ws.Range("M73").Value = -3512.3335

$ws.Range("H80").Value = 4029.1667
$ws.Range("I80").Value = 3694.7368
$ws.Range("J80").Value = 5300
$ws.Range("K80").Value = 3694.7368
$ws.Range("L80").Value = 5300
$ws.Range("M80").Value = -2696.7368
$ws.Range("N80").Value = -7296

$ws.Range("H83").Value = 4029.1667
$ws.Range("I83").Value = 3694.7368
$ws.Range("J83").Value = 5300
$ws.Range("K83").Value = 18473.684
$ws.Range("L83").Value = 26500
$ws.Range("M83").Value = -13481.684
$ws.Range("N83").Value = -36484

$ws.Range("H102").Value = 1326.7333
$ws.Range("I102").Value = 1386
$ws.Range("J102").Value = 1274.875
$ws.Range("K102").Value = 1386
$ws.Range("L102").Value = 1274.875
$ws.Range("M102").Value = 236
$ws.Range("N102").Value = -4518.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5961.125
$ws.Range("I16").Value = 7382.4116
$ws.Range("J16").Value = 2509.4285
$ws.Range("K16").Value = 7382.4116
$ws.Range("L16").Value = 2509.4285
$ws.Range("M16").Value = -7212.4116
$ws.Range("N16").Value = -2849.4285

$ws.Range("H55").Value = 192
$ws.Range("I55").Value = 86.5
$ws.Range("J55").Value = 271.125
$ws.Range("K55").Value = 86.5
$ws.Range("L55").Value = 271.125
$ws.Range("M55").Value = 86.5
$ws.Range("N55").Value = -617.125

$ws.Range("H68").Value = 4212.5
$ws.Range("I68").Value = 4475.5
$ws.Range("J68").Value = 3949.5
$ws.Range("K68").Value = 4475.5
$ws.Range("L68").Value = 3949.5
$ws.Range("M68").Value = -3726.5
$ws.Range("N68").Value = -5447.5

$ws.Range("H71").Value = 4212.5
$ws.Range("I71").Value = 4475.5
$ws.Range("J71").Value = 3949.5
$ws.Range("K71").Value = 22377.5
$ws.Range("L71").Value = 19747.5
$ws.Range("M71").Value = -18633.5
$ws.Range("N71").Value = -27235.5

$ws.Range("H132").Value = 1701.6271
$ws.Range("I132").Value = 1467.96
$ws.Range("J132").Value = 2999.7778
$ws.Range("K132").Value = 4403.88
$ws.Range("L132").Value = 8999.3334
$ws.Range("M132").Value = -1873.88
$ws.Range("N132").Value = -14059.3334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2665
$ws.Range("I81").Value = 2596
$ws.Range("J81").Value = 2714.2856
$ws.Range("K81").Value = 5192
$ws.Range("L81").Value = 5428.5712
$ws.Range("M81").Value = -4131
$ws.Range("N81").Value = -7550.5712

$ws.Range("H84").Value = 2665
$ws.Range("I84").Value = 2596
$ws.Range("J84").Value = 2714.2856
$ws.Range("K84").Value = 25960
$ws.Range("L84").Value = 27142.856
$ws.Range("M84").Value = -20656
$ws.Range("N84").Value = -37750.856

$ws.Range("H86").Value = 16331.25
$ws.Range("J86").Value = 16331.25
$ws.Range("L86").Value = 16331.25
$ws.Range("N86").Value = -18577.25

$ws.Range("H89").Value = 16331.25
$ws.Range("J89").Value = 16331.25
$ws.Range("L89").Value = 81656.25
$ws.Range("N89").Value = -92888.25

$ws.Range("H132").Value = 15874346
$ws.Range("I132").Value = 33334374
$ws.Range("J132").Value = 1593.4849
$ws.Range("K132").Value = 100003122
$ws.Range("L132").Value = 4780.4547
$ws.Range("M132").Value = -100000592
$ws.Range("N132").Value = -9840.4547
